$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.626.83"
$ws.Range("E2").Value = "  -0.17%  "

$ws.Range("D3").Value = "1.592.23"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("E9").Value = "  -1.87%  "

$ws.Range("E10").Value = "  -1.31%  "

$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "1.816.15"

$ws.Range("D13").Value = "1.575.17"
$ws.Range("E13").Value = "  -0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "26.604.12"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "206.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.95%  "

$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("E26").Value = "  +0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.661"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").Value = "1.277.90"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0172"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.919"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.94%  "

$ws.Range("D46").Value = "1.728.70"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "

$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.96%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
